# Apply the "closer to finish tests" edit to sheet1:
#  - rewrite column A test-case names (f1..f4 -> grid/timeSeries/timeSeriesProfile/trajectory)
#  - insert a new "fail" row after the four feature-type rows
#  - relabel the "basiert auf" column ("f1" -> "grid", or "-" where no longer applicable)
#  - widen column A, move the active selection, matching the new used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite rows 2-5 (feature type test cases) ---
$ws.Range("A2").Value = "grid"

$ws.Range("A3").Value = "timeSeries"

$ws.Range("A4").Value = "timeSeriesProfile"

$ws.Range("A5").Value = "trajectory"

# --- Insert a new row 6 for the "fail" test case, shifting old rows 6-11 to 7-12 ---
[void]$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = "fail"
$ws.Range("B3").Value = "-"
$ws.Range("D6").Value = "alles falsch"

$ws.Range("B4").Value = "-"
$ws.Range("B5").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"

# --- Fix up "basiert auf" column for the rows that used to read "f1" ---
$ws.Range("B7").Value = "grid"
$ws.Range("B8").Value = "grid"
$ws.Range("B12").Value = "grid"

# --- Column A width / bestFit ---
$ws.Columns("A:A").ColumnWidth = 22.85

# --- Selection moves to B13 ---
[void]$ws.Range("B13").Select()
